$d = $word.ActiveDocument

$d.Content.Find.Execute("Fornecer conhecimentos que proporcionam uma visão holística a respeito da Cadeia de Suprimentos e da Logística, apresentando métodos e ferramentas para otimizar o desempenho das cadeias produtivas.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_0@@", 2) | Out-Null
$d.Content.Find.Execute("Provide knowledge that paves a holistic view of Logistics and Supply Chain, presenting methods and tools available to optimize the performance of production chains.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_1@@", 2) | Out-Null
$d.Content.Find.Execute("3295113 - José Eduardo Holler Branco", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_2@@", 2) | Out-Null
$d.Content.Find.Execute("Gerenciamento da cadeia de suprimentos e da logística: planejamento, otimização e controle.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_3@@", 2) | Out-Null
$d.Content.Find.Execute("Supply chain and logistics management: planning, optimization and control.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_4@@", 2) | Out-Null
$d.Content.Find.Execute("i) Introdução à Logística e Cadeia de Suprimentos; ii) Planejamento da cadeia de suprimentos; iii) Planejamento do transporte; iv) Custos logísticos; v) Tipos de cargas e sistemas de armazenamento; vi) Modelos de transporte; vii) Modelos de localização; viii) Planejamento do estoque; ix) Logística Reversa e Economia Circular; e x)  Controle da logística e cadeia de suprimentos.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_5@@", 2) | Out-Null
$d.Content.Find.Execute("Provas, trabalhos em grupo, exercícios individuais e seminários.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_6@@", 2) | Out-Null
$d.Content.Find.Execute("Média das atividades avaliativas.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_7@@", 2) | Out-Null
$d.Content.Find.Execute("NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_8@@", 2) | Out-Null
$d.Content.Find.Execute("BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gestão Logística da Cadeia de Suprimentos. 4. ed. AMGH, 2013. 472 p.", $true, $false, $false, $false, $false, $true, 1, $false, "@@SWAP_MARK_9@@", 2) | Out-Null

$d.Content.Find.Execute("@@SWAP_MARK_0@@", $true, $false, $false, $false, $false, $true, 1, $false, "Gerenciamento da cadeia de suprimentos e da logística: planejamento, otimização e controle.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_1@@", $true, $false, $false, $false, $false, $true, 1, $false, "Supply chain and logistics management: planning, optimization and control.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_2@@", $true, $false, $false, $false, $false, $true, 1, $false, "Fornecer conhecimentos que proporcionam uma visão holística a respeito da Cadeia de Suprimentos e da Logística, apresentando métodos e ferramentas para otimizar o desempenho das cadeias produtivas.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_3@@", $true, $false, $false, $false, $false, $true, 1, $false, "i) Introdução à Logística e Cadeia de Suprimentos; ii) Planejamento da cadeia de suprimentos; iii) Planejamento do transporte; iv) Custos logísticos; v) Tipos de cargas e sistemas de armazenamento; vi) Modelos de transporte; vii) Modelos de localização; viii) Planejamento do estoque; ix) Logística Reversa e Economia Circular; e x)  Controle da logística e cadeia de suprimentos.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_4@@", $true, $false, $false, $false, $false, $true, 1, $false, "Provide knowledge that paves a holistic view of Logistics and Supply Chain, presenting methods and tools available to optimize the performance of production chains.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_5@@", $true, $false, $false, $false, $false, $true, 1, $false, "Provas, trabalhos em grupo, exercícios individuais e seminários.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_6@@", $true, $false, $false, $false, $false, $true, 1, $false, "Média das atividades avaliativas.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_7@@", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_8@@", $true, $false, $false, $false, $false, $true, 1, $false, "BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gestão Logística da Cadeia de Suprimentos. 4. ed. AMGH, 2013. 472 p.", 2) | Out-Null
$d.Content.Find.Execute("@@SWAP_MARK_9@@", $true, $false, $false, $false, $false, $true, 1, $false, "3295113 - José Eduardo Holler Branco", 2) | Out-Null

Write-Output "done"